# "Some Changes in Workflow"
#
# - Rename the "Image Link" header (E1) to "ImageLink"
# - Refresh several postimg.cc map-image URLs in column E (rows 2, 4, 7, 9)
# - Refresh a couple of "Hours" strings (column H, rows 4 and 7) to the
#   "Closes soon" phrasing
# - Remove the hidden legacy "Packager Shell Object" OLE icon (anchored at I2)
#   together with its drawing
# - Land the selection back on F6 (and drop the custom topLeftCell scroll
#   position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename -------------------------------------------------------
$ws.Range("E1").Value = "ImageLink"

# --- Updated map image links --------------------------------------------
$ws.Range("E2").Value = "https://i.postimg.cc/9fwWRQDD/1map.png"
$ws.Range("E4").Value = "https://i.postimg.cc/sfZDqtMy/3map.png"
$ws.Range("E7").Value = "https://i.postimg.cc/V6W4k76N/6map.png"
$ws.Range("E9").Value = "https://i.postimg.cc/505xxKry/8map.png"

# --- Updated "Hours" text -------------------------------------------------
$ws.Range("H4").Value = "       Closes soon ⋅ 6PM ⋅ Opens 10AM Wed        "
$ws.Range("H7").Value = "       Closes soon ⋅ 6:30PM ⋅ Opens 10AM Wed        "

# --- Drop the hidden embedded OLE "Packager Shell Object" icon ----------
if ($ws.Shapes.Count -gt 0) {
    for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
        $ws.Shapes.Item($i).Delete()
    }
}

# --- Move the active selection to F6 and reset the top-left scroll cell -
[void]$ws.Range("F6").Select()
